# Add "arrival" / "departure" date columns (P, Q) to the activity sheet,
# each formatted as dd/mm/yy, and fill in the corresponding dates for the
# existing data rows. Also drops the stray fully-empty last row
# (row 1048576) that was lingering in the sheet, and leaves the cursor
# on A6 as the new active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the stray trailing empty row that has no data.
$ws.Rows.Item(1048576).Delete()

# New header cells.
$ws.Range("P1").Value = "arrival"
$ws.Range("Q1").Value = "departure"

# New dated cells - set the number format first so the value is stored
# using that format straight away (dd/mm/yy), then assign the date.
$ws.Range("P2").NumberFormat = "dd/mm/yy"
$ws.Range("P2").Value = [DateTime]"2000-01-01"

$ws.Range("P3").NumberFormat = "dd/mm/yy"
$ws.Range("P3").Value = [DateTime]"2010-02-01"

$ws.Range("Q4").NumberFormat = "dd/mm/yy"
$ws.Range("Q4").Value = [DateTime]"2019-05-12"

$ws.Range("P5").NumberFormat = "dd/mm/yy"
$ws.Range("P5").Value = [DateTime]"2012-01-01"

$ws.Range("Q6").NumberFormat = "dd/mm/yy"
$ws.Range("Q6").Value = [DateTime]"2019-05-12"

$ws.Range("Q7").NumberFormat = "dd/mm/yy"
$ws.Range("Q7").Value = [DateTime]"2019-05-12"

# Leave the selection where the author left it.
$ws.Range("A6").Select() | Out-Null
